$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) and Volume (E) columns for rows 2-51 to avoid Excel
# auto-converting numeric-looking strings (prices, percentages) into actual numbers,
# which would destroy formatting like "3.068.57", "0.0000229", leading zeros, "0.0₃0520", etc.
$ws.Range("D2:E51").NumberFormat = "@"


# Row 2
$ws.Range("D2").Value = '64.220.58'
$ws.Range("E2").Value = '  +1.19%  '

# Row 3
$ws.Range("D3").Value = '3.073.99'
$ws.Range("E3").Value = '  +0.38%  '

# Row 4
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.13%  '

# Row 5
$ws.Range("D5").Value = '557.79'
$ws.Range("E5").Value = '  +1.62%  '

# Row 6
$ws.Range("D6").Value = '146.29'
$ws.Range("E6").Value = '  +5.31%  '

# Row 7
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.22%  '

# Row 8
$ws.Range("D8").Value = '3.075.05'
$ws.Range("E8").Value = '  +0.65%  '

# Row 9
$ws.Range("E9").Value = '  +0.54%  '

# Row 10
$ws.Range("B10").Value = 'Toncoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D10").Value = '6.34'
$ws.Range("E10").Value = '  +1.70%  '

# Row 11
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '0.153'
$ws.Range("E11").Value = '  +2.56%  '

# Row 12
$ws.Range("E12").Value = '  +3.80%  '

# Row 13
$ws.Range("D13").Value = '0.0000229'
$ws.Range("E13").Value = '  +0.60%  '

# Row 14
$ws.Range("D14").Value = '35.32'
$ws.Range("E14").Value = '  +1.89%  '

# Row 15
$ws.Range("D15").Value = '3.570.68'
$ws.Range("E15").Value = '  +0.12%  '

# Row 16
$ws.Range("D16").Value = '64.133.62'
$ws.Range("E16").Value = '  +0.94%  '

# Row 17
$ws.Range("D17").Value = '3.066.07'
$ws.Range("E17").Value = '  -0.07%  '

# Row 18
$ws.Range("E18").Value = '  +1.22%  '

# Row 19
$ws.Range("D19").Value = '6.80'
$ws.Range("E19").Value = '  +1.16%  '

# Row 20
$ws.Range("D20").Value = '478.33'
$ws.Range("E20").Value = '  -0.14%  '

# Row 21
$ws.Range("D21").Value = '13.97'
$ws.Range("E21").Value = '  +3.08%  '

# Row 22
$ws.Range("D22").Value = '0.677'
$ws.Range("E22").Value = '  -0.27%  '

# Row 23
$ws.Range("D23").Value = '7.56'
$ws.Range("E23").Value = '  +4.99%  '

# Row 24
$ws.Range("D24").Value = '13.56'
$ws.Range("E24").Value = '  +8.35%  '

# Row 25
$ws.Range("D25").Value = '81.70'
$ws.Range("E25").Value = '  +0.34%  '

# Row 26
$ws.Range("E26").Value = '  +0.15%  '

# Row 27
$ws.Range("D27").Value = '2.80'
$ws.Range("E27").Value = '  +1.67%  '

# Row 28
$ws.Range("D28").Value = '8.11'
$ws.Range("E28").Value = '  +2.02%  '

# Row 29
$ws.Range("D29").Value = '2.07'
$ws.Range("E29").Value = '  +4.27%  '

# Row 30
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.11%  '

# Row 31
$ws.Range("D31").Value = '26.24'
$ws.Range("E31").Value = '  +1.04%  '

# Row 32
$ws.Range("E32").Value = '  +0.87%  '

# Row 33
$ws.Range("D33").Value = '2.50'
$ws.Range("E33").Value = '  +3.46%  '

# Row 34
$ws.Range("D34").Value = '5.60'
$ws.Range("E34").Value = '  -1.22%  '

# Row 35
$ws.Range("D35").Value = '6.19'
$ws.Range("E35").Value = '  +3.40%  '

# Row 36
$ws.Range("D36").Value = '54.91'
$ws.Range("E36").Value = '  -1.56%  '

# Row 37
$ws.Range("D37").Value = '462.58'
$ws.Range("E37").Value = '  -1.07%  '

# Row 38
$ws.Range("D38").Value = '3.04'
$ws.Range("E38").Value = '  +18.04%  '

# Row 39
$ws.Range("D39").Value = '0.0834'
$ws.Range("E39").Value = '  +2.28%  '

# Row 40
$ws.Range("D40").Value = '0.0406'
$ws.Range("E40").Value = '  +2.74%  '

# Row 41
$ws.Range("D41").Value = '2.964.87'
$ws.Range("E41").Value = '  -5.63%  '

# Row 42
$ws.Range("D42").Value = '8.30'
$ws.Range("E42").Value = '  +0.78%  '

# Row 43
$ws.Range("D43").Value = '0.115'
$ws.Range("E43").Value = '  -3.81%  '

# Row 44
$ws.Range("D44").Value = '28.21'
$ws.Range("E44").Value = '  +0.91%  '

# Row 45
$ws.Range("D45").Value = '0.262'
$ws.Range("E45").Value = '  +4.29%  '

# Row 46
$ws.Range("D46").Value = '2.16'
$ws.Range("E46").Value = '  +5.20%  '

# Row 47
$ws.Range("E47").Value = '  +0.04%  '

# Row 48
$ws.Range("E48").Value = '  +2.76%  '

# Row 49
$ws.Range("D49").Value = '120.06'
$ws.Range("E49").Value = '  +3.29%  '

# Row 50
$ws.Range("D50").Value = ('0.0' + [string][char]0x2083 + '0520')
$ws.Range("E50").Value = '  +1.69%  '

# Row 51
$ws.Range("D51").Value = '2.09'
$ws.Range("E51").Value = '  +1.12%  '
